$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step date bump (A1)
$ws.Range("A1").Value = 45309

# Updated quantity/price values in column D (rows 14-26, 28-31)
$ws.Range("D14").Value = 7835.894
$ws.Range("D15").Value = 8730.41
$ws.Range("D16").Value = 10891.583
$ws.Range("D17").Value = 13238.733
$ws.Range("D18").Value = 15224.515
$ws.Range("D19").Value = 16119.026
$ws.Range("D20").Value = 23277.144
$ws.Range("D21").Value = 20036.993
$ws.Range("D22").Value = 23507.7
$ws.Range("D23").Value = 25361.078
$ws.Range("D24").Value = 18736.258
$ws.Range("D25").Value = 17442.959
$ws.Range("D26").Value = 20466.318
$ws.Range("D28").Value = 20931.508
$ws.Range("D29").Value = 25260.957
$ws.Range("D30").Value = 31343.565
$ws.Range("D31").Value = 36567.549
